# Normalize the "Recorded By" (column G) lists so that any "System" /
# "system" token is moved to the end of the list, and the remaining
# (non-system) tokens are ordered in descending alphabetical order.
#
# Example:
#   "System, backup@backdoor.com, system" -> "backup@backdoor.com, System, system"
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"

function Transform-RecordedBy {
    param([string]$val)

    if ([string]::IsNullOrEmpty($val)) {
        return $val
    }

    $parts = @($val.Split(","))
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $nonSystem = @()
    $systemItems = @()
    foreach ($t in $trimmed) {
        if ($t.ToLower() -eq "system") {
            $systemItems += $t
        } else {
            $nonSystem += $t
        }
    }

    $sortedNonSystem = @($nonSystem | Sort-Object -Descending)
    $result = @($sortedNonSystem) + @($systemItems)

    return ($result -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value2
    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    $new = Transform-RecordedBy $current
    if ($new -ne $current) {
        $cell.Value = $new
    }
}
